$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '90.037.64'
Set-TextValue 'E2' '  +2.96%  '
Set-TextValue 'D3' '3.201.09'
Set-TextValue 'E3' '  -0.73%  '
Set-TextValue 'E4' '  -0.03%  '
Set-TextValue 'D5' '217.48'
Set-TextValue 'E5' '  +6.03%  '
Set-TextValue 'D6' '623.59'
Set-TextValue 'E6' '  +2.04%  '
Set-TextValue 'D7' '0.388'
Set-TextValue 'E7' '  +2.48%  '
Set-TextValue 'D8' '0.695'
Set-TextValue 'E8' '  +3.08%  '
Set-TextValue 'D9' '0.999'
Set-TextValue 'E9' '  +0.06%  '
Set-TextValue 'D10' '3.196.85'
Set-TextValue 'E10' '  -0.75%  '
Set-TextValue 'E11' '  +5.39%  '
Set-TextValue 'E12' '  -0.58%  '
Set-TextValue 'D13' '0.0000254'
Set-TextValue 'E13' '  +2.84%  '
Set-TextValue 'D14' '5.39'
Set-TextValue 'E14' '  +1.45%  '
Set-TextValue 'B15' 'WrappedBTC'
Set-TextValue 'C15' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D15' '89.830.00'
Set-TextValue 'E15' '  +2.85%  '
Set-TextValue 'B16' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C16' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D16' '3.793.68'
Set-TextValue 'E16' '  -0.61%  '
Set-TextValue 'D17' '32.91'
Set-TextValue 'E17' '  +0.76%  '
Set-TextValue 'D18' '3.210.97'
Set-TextValue 'E18' '  -0.84%  '
Set-TextValue 'D19' '0.0000236'
Set-TextValue 'E19' '  +77.04%  '
Set-TextValue 'D20' '3.38'
Set-TextValue 'E20' '  +14.22%  '
Set-TextValue 'D21' '13.39'
Set-TextValue 'E21' '  -0.58%  '
Set-TextValue 'D22' '436.28'
Set-TextValue 'E22' '  +3.60%  '
Set-TextValue 'D23' '8.55'
Set-TextValue 'E23' '  -0.34%  '
Set-TextValue 'D24' '5.06'
Set-TextValue 'E24' '  -1.41%  '
Set-TextValue 'D25' '5.10'
Set-TextValue 'E25' '  -2.10%  '
Set-TextValue 'D26' '11.62'
Set-TextValue 'E26' '  -0.19%  '
Set-TextValue 'B27' 'Litecoin'
Set-TextValue 'C27' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D27' '78.13'
Set-TextValue 'E27' '  +5.08%  '
Set-TextValue 'B28' 'WrappedeETH'
Set-TextValue 'C28' 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue 'D28' '3.356.90'
Set-TextValue 'E28' '  -1.01%  '
Set-TextValue 'E29' '  +0.33%  '
Set-TextValue 'D30' '0.998'
Set-TextValue 'E30' '  -0.34%  '
Set-TextValue 'D31' '0.157'
Set-TextValue 'E31' '  -9.18%  '
Set-TextValue 'D32' '4.09'
Set-TextValue 'E32' '  +36.08%  '
Set-TextValue 'D33' '8.45'
Set-TextValue 'E33' '  +0.33%  '
Set-TextValue 'D34' '536.11'
Set-TextValue 'E34' '  -1.66%  '
Set-TextValue 'D35' '1.89'
Set-TextValue 'E35' '  +0.64%  '
Set-TextValue 'D36' '6.86'
Set-TextValue 'E36' '  +2.89%  '
Set-TextValue 'D37' '1.28'
Set-TextValue 'E37' '  -0.60%  '
Set-TextValue 'D38' '22.34'
Set-TextValue 'E38' '  +0.34%  '
Set-TextValue 'E39' '  +2.17%  '
Set-TextValue 'D40' '0.999'
Set-TextValue 'E40' '  +0.09%  '
Set-TextValue 'D41' '0.127'
Set-TextValue 'E41' '  -6.08%  '
Set-TextValue 'E42' '  +0.05%  '
Set-TextValue 'D43' '1.92'
Set-TextValue 'E43' '  +1.10%  '
Set-TextValue 'D44' '0.371'
Set-TextValue 'E44' '  -2.83%  '
Set-TextValue 'D45' '149.72'
Set-TextValue 'E45' '  +2.18%  '
Set-TextValue 'D46' '171.61'
Set-TextValue 'E46' '  -1.72%  '
Set-TextValue 'D47' '43.27'
Set-TextValue 'E47' '  -0.68%  '
Set-TextValue 'E48' '  -4.10%  '
Set-TextValue 'D49' '1.24'
Set-TextValue 'E49' '  -3.51%  '
Set-TextValue 'D50' '0.734'
Set-TextValue 'E50' '  +4.20%  '
Set-TextValue 'D51' '0.614'
Set-TextValue 'E51' '  +1.62%  '
